# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table: Good Roaming Calculation (%) tweaks ---
$ws.Range("D3").Value = 83.8
$ws.Range("D5").Value = 95.7

# --- Good Drivers table: refreshed driver list (rows 14-22) ---
# A new driver entry lands on top (row 14); everything else rotates down one
# row, and three rows that previously had no "Driver Vintage" date now get
# one. Column E holds literal text dates, so format it as Text first so the
# COM layer doesn't silently reinterpret the strings as date serials.
$ws.Range("E14:E21").NumberFormat = "@"

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B14").Value = 445055
$ws.Range("D14").Value = 99.90000000000001
$ws.Range("E14").Value = "2024-11-10"

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B15").Value = 10661
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = "2022-08-29"

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B16").Value = 14239
$ws.Range("D16").Value = 100
$ws.Range("E16").Value = "2022-05-23"

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B17").Value = 265400
$ws.Range("D17").Value = 99.90000000000001
$ws.Range("E17").Value = "2022-05-01"

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B18").Value = 77849
$ws.Range("D18").Value = 99.90000000000001
$ws.Range("E18").Value = "2021-08-18"

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B19").Value = 34244
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = "2021-04-27"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B20").Value = 59673
$ws.Range("D20").Value = 100
$ws.Range("E20").Value = "2020-08-05"

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B21").Value = 113652
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = "2020-01-06"

$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B22").Value = 56018
$ws.Range("D22").Value = 100
# E22 ("2019-12-14") is unchanged from the original report.
